$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.5325700330326129
$ws.Range("C2").Value = 0.07171134690607062
$ws.Range("E2").Value = 0.09701186566914544
$ws.Range("F2").Value = 0.4443680307746263
$ws.Range("G2").Value = 1.376991519755975
$ws.Range("H2").Value = 1.300496466101066
$ws.Range("I2").Value = 1.325132537853463
$ws.Range("K2").Value = 0.3345138496387108
$ws.Range("L2").Value = 0.2184663207233086
$ws.Range("M2").Value = 0.164825157848302
$ws.Range("B3").Value = 0.5009736353710537
$ws.Range("C3").Value = 0.06909471166952841
$ws.Range("E3").Value = 0.09723511531864659
$ws.Range("F3").Value = 0.387822817061874
$ws.Range("G3").Value = 1.384807818900853
$ws.Range("H3").Value = 1.309205032845469
$ws.Range("I3").Value = 1.335038009770244
$ws.Range("K3").Value = 0.3035680431421497
$ws.Range("L3").Value = 0.215948427421857
$ws.Range("M3").Value = 0.1587435054090882
$ws.Range("B4").Value = 0.4817774397117489
$ws.Range("C4").Value = 0.06746485726785068
$ws.Range("E4").Value = 0.09740470357287734
$ws.Range("F4").Value = 0.3531389305168915
$ws.Range("G4").Value = 1.390361661925581
$ws.Range("H4").Value = 1.315073714672465
$ws.Range("I4").Value = 1.341709805329064
$ws.Range("K4").Value = 0.2846535573075499
$ws.Range("L4").Value = 0.2145079094339408
$ws.Range("M4").Value = 0.1550847539822477
$ws.Range("B5").Value = 0.4740066255261866
$ws.Range("C5").Value = 0.06679482544931403
$ws.Range("E5").Value = 0.09748200461575429
$ws.Range("F5").Value = 0.3390132514313251
$ws.Range("G5").Value = 1.392814494098246
$ws.Range("H5").Value = 1.317596442192482
$ws.Range("I5").Value = 1.344576903353701
$ws.Range("K5").Value = 0.2769677922792795
$ws.Range("L5").Value = 0.2139474631806735
$ws.Range("M5").Value = 0.1536128385240936
$ws.Range("B6").Value = 0.47271942954103
$ws.Range("C6").Value = 0.06668321334144167
$ws.Range("E6").Value = 0.0974953357002164
$ws.Range("F6").Value = 0.336668177824194
$ws.Range("G6").Value = 1.393233231608662
$ws.Range("H6").Value = 1.318023263887738
$ws.Range("I6").Value = 1.345061938349119
$ws.Range("K6").Value = 0.2756929177223668
$ws.Range("L6").Value = 0.2138560082309908
$ws.Range("M6").Value = 0.1533695817527203
$ws.Range("B7").Value = 0.4816724294749406
$ws.Range("C7").Value = 0.06745584469685184
$ws.Range("E7").Value = 0.09740571288795508
$ws.Range("F7").Value = 0.3529483938344953
$ws.Range("G7").Value = 1.390393974219378
$ws.Range("H7").Value = 1.315107205823622
$ws.Range("I7").Value = 1.341747871638013
$ws.Range("K7").Value = 0.2845498147147509
$ws.Range("L7").Value = 0.2145002433813517
$ws.Range("M7").Value = 0.1550648259361616
$ws.Range("B8").Value = 0.5216334965336387
$ws.Range("C8").Value = 0.0708139447546543
$ws.Range("E8").Value = 0.09708210436678932
$ws.Range("F8").Value = 0.4248636149813336
$ws.Range("G8").Value = 1.379529867056149
$ws.Range("H8").Value = 1.303390945210637
$ws.Range("I8").Value = 1.328425519211031
$ws.Range("K8").Value = 0.3238259946662367
$ws.Range("L8").Value = 0.2175762894750974
$ws.Range("M8").Value = 0.1627126069354006
$ws.Range("B9").Value = 0.6016007230618925
$ws.Range("C9").Value = 0.07721564099557554
$ws.Range("E9").Value = 0.09670477869301486
$ws.Range("F9").Value = 0.5661985755041457
$ws.Range("G9").Value = 1.364220585244894
$ws.Range("H9").Value = 1.284552643516648
$ws.Range("I9").Value = 1.306981394634171
$ws.Range("K9").Value = 0.4015215805349328
$ws.Range("L9").Value = 0.2244435680668531
$ws.Range("M9").Value = 0.178305041710594
$ws.Range("B10").Value = 0.6613155085839537
$ws.Range("C10").Value = 0.08180856166585215
$ws.Range("E10").Value = 0.0965835288228476
$ws.Range("F10").Value = 0.6702781546542269
$ws.Range("G10").Value = 1.356639524202023
$ws.Range("H10").Value = 1.273232903212673
$ws.Range("I10").Value = 1.294082031993348
$ws.Range("K10").Value = 0.4590085242639077
$ws.Range("L10").Value = 0.2299965601603446
$ws.Range("M10").Value = 0.1901207641981486
$ws.Range("B11").Value = 0.6886876983664365
$ws.Range("C11").Value = 0.08387436762549783
$ws.Range("E11").Value = 0.09656205866543743
$ws.Range("F11").Value = 0.7176906081379002
$ws.Range("G11").Value = 1.353989515841718
$ws.Range("H11").Value = 1.268630369126981
$ws.Range("I11").Value = 1.288834327234902
$ws.Range("K11").Value = 0.4852473483947506
$ws.Range("L11").Value = 0.2326327117797291
$ws.Range("M11").Value = 0.1955736905285903
$ws.Range("B12").Value = 0.6990823250916662
$ws.Range("C12").Value = 0.08465326862983602
$ws.Range("E12").Value = 0.09655875678478765
$ws.Range("F12").Value = 0.7356546913071611
$ws.Range("G12").Value = 1.353101071585371
$ws.Range("H12").Value = 1.266966131508042
$ws.Range("I12").Value = 1.286936396056419
$ws.Range("K12").Value = 0.4951956889036921
$ws.Range("L12").Value = 0.2336467414921657
$ws.Range("M12").Value = 0.1976497001880944
$ws.Range("B13").Value = 0.6968423575866325
$ws.Range("C13").Value = 0.08448566830638526
$ws.Range("E13").Value = 0.09655925340987892
$ws.Range("F13").Value = 0.7317853510981394
$ws.Range("G13").Value = 1.353287293432885
$ws.Range("H13").Value = 1.267321057177256
$ws.Range("I13").Value = 1.287341178557377
$ws.Range("K13").Value = 0.4930525956965255
$ws.Range("L13").Value = 0.233427651210647
$ws.Range("M13").Value = 0.1972021020044536
$ws.Range("B14").Value = 0.6895422846631618
$ws.Range("C14").Value = 0.08393851586160395
$ws.Range("E14").Value = 0.09656169033940998
$ws.Range("F14").Value = 0.7191683204515869
$ws.Range("G14").Value = 1.353914115889737
$ws.Range("H14").Value = 1.268491875173922
$ws.Range("I14").Value = 1.288676394002302
$ws.Range("K14").Value = 0.4860655604841213
$ws.Range("L14").Value = 0.2327158207268241
$ws.Range("M14").Value = 0.1957442630864676
$ws.Range("B15").Value = 0.6850745900310358
$ws.Range("C15").Value = 0.08360292997035401
$ws.Range("E15").Value = 0.09656381136295877
$ws.Range("F15").Value = 0.7114413442032514
$ws.Range("G15").Value = 1.354313052735847
$ws.Range("H15").Value = 1.269219276619026
$ws.Range("I15").Value = 1.289505878337025
$ws.Range("K15").Value = 0.4817873869544655
$ws.Range("L15").Value = 0.2322818570224996
$ws.Range("M15").Value = 0.1948527379748413
$ws.Range("B16").Value = 0.6595308082148676
$ws.Range("C16").Value = 0.08167308359568892
$ws.Range("E16").Value = 0.0965856084697414
$ws.Range("F16").Value = 0.6671810134426437
$ws.Range("G16").Value = 1.356828797480617
$ws.Range("H16").Value = 1.2735446943162
$ws.Range("I16").Value = 1.294437469028111
$ws.Range("K16").Value = 0.4572954908980478
$ws.Range("L16").Value = 0.2298264919081277
$ws.Range("M16").Value = 0.1897659618332099
$ws.Range("B17").Value = 0.6439133595579563
$ws.Range("C17").Value = 0.08048316388797616
$ws.Range("E17").Value = 0.09660759712063971
$ws.Range("F17").Value = 0.6400460337125793
$ws.Range("G17").Value = 1.358576830077197
$ws.Range("H17").Value = 1.27633826676329
$ws.Range("I17").Value = 1.297621755085075
$ws.Range("K17").Value = 0.4422927265198666
$ws.Range("L17").Value = 0.2283483593826503
$ws.Range("M17").Value = 0.1866652714607326
$ws.Range("B18").Value = 0.6349501885030122
$ws.Range("C18").Value = 0.07979653757699623
$ws.Range("E18").Value = 0.09662341613186953
$ws.Range("F18").Value = 0.6244449056556647
$ws.Range("G18").Value = 1.359657412722129
$ws.Range("H18").Value = 1.277996529778363
$ws.Range("I18").Value = 1.299511651967336
$ws.Range("K18").Value = 0.4336718109938431
$ws.Range("L18").Value = 0.2275085397996151
$ws.Range("M18").Value = 0.1848891739405616
$ws.Range("B19").Value = 0.6319187891872105
$ws.Range("C19").Value = 0.0795636768177701
$ws.Range("E19").Value = 0.09662931740375491
$ws.Range("F19").Value = 0.619163680173358
$ws.Range("G19").Value = 1.360036181906437
$ws.Range("H19").Value = 1.278566829793448
$ws.Range("I19").Value = 1.300161563381174
$ws.Range("K19").Value = 0.4307543494572883
$ws.Range("L19").Value = 0.227225972764316
$ws.Range("M19").Value = 0.1842890810170914
$ws.Range("B20").Value = 0.6455738414871064
$ws.Range("C20").Value = 0.08061006214711597
$ws.Range("E20").Value = 0.09660492824450451
$ws.Range("F20").Value = 0.642933953830422
$ws.Range("G20").Value = 1.358382968658773
$ws.Range("H20").Value = 1.276035558818919
$ws.Range("I20").Value = 1.297276739837457
$ws.Range("K20").Value = 0.4438889422095826
$ws.Range("L20").Value = 0.2285046369384958
$ws.Range("M20").Value = 0.1869945865329399
$ws.Range("B21").Value = 0.6916856997277421
$ws.Range("C21").Value = 0.08409931931350911
$ws.Range("E21").Value = 0.09656084363791528
$ws.Range("F21").Value = 0.7228739723491628
$ws.Range("G21").Value = 1.353726878662854
$ws.Range("H21").Value = 1.268145843198354
$ws.Range("I21").Value = 1.288281785828282
$ws.Range("K21").Value = 0.4881174919225657
$ws.Range("L21").Value = 0.2329244748387396
$ws.Range("M21").Value = 0.1961721651000872
$ws.Range("B22").Value = 0.7219934470418536
$ws.Range("C22").Value = 0.0863600837735703
$ws.Range("E22").Value = 0.09656016814499679
$ws.Range("F22").Value = 0.7751780083420101
$ws.Range("G22").Value = 1.351354589081808
$ws.Range("H22").Value = 1.263447846238222
$ws.Range("I22").Value = 1.282923373690991
$ws.Range("K22").Value = 0.5170947492226219
$ws.Range("L22").Value = 0.2359050270650869
$ws.Range("M22").Value = 0.2022349260969563
$ws.Range("B23").Value = 0.7058021465651052
$ws.Range("C23").Value = 0.08515526752482572
$ws.Range("E23").Value = 0.09655795938264156
$ws.Range("F23").Value = 0.7472568307830727
$ws.Range("G23").Value = 1.352559282841099
$ws.Range("H23").Value = 1.265913312068577
$ws.Range("I23").Value = 1.285735629046258
$ws.Range("K23").Value = 0.5016226279860234
$ws.Range("L23").Value = 0.2343058553237398
$ws.Range("M23").Value = 0.1989932296991057
$ws.Range("B24").Value = 0.6448230892653157
$ws.Range("C24").Value = 0.08055269931372777
$ws.Range("E24").Value = 0.09660612494495702
$ws.Range("F24").Value = 0.6416283278902171
$ws.Range("G24").Value = 1.358470377880508
$ws.Range("H24").Value = 1.276172250473167
$ws.Range("I24").Value = 1.297432536801111
$ws.Range("K24").Value = 0.4431672794484882
$ws.Range("L24").Value = 0.2284339527704162
$ws.Range("M24").Value = 0.1868456828395111
$ws.Range("B25").Value = 0.5797973232353115
$ws.Range("C25").Value = 0.07550330810440187
$ws.Range("E25").Value = 0.0967794026278419
$ws.Range("F25").Value = 0.5279251897347166
$ws.Range("G25").Value = 1.367719100514933
$ws.Range("H25").Value = 1.289206131550287
$ws.Range("I25").Value = 1.312281264266502
$ws.Range("K25").Value = 0.380431381836047
$ws.Range("L25").Value = 0.2224965365962817
$ws.Range("M25").Value = 0.1740234570257435
